$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5 and 6 (FABRICIO / 004570632 and MARCIA / 004547722)
$ws.Rows("5:6").Delete()
